# Update the "Förändrad" (Changed) date column (C) for rows 2-101
# from 2023-10-05 (serial 45204) to 2023-10-08 (serial 45207).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSerial = 45204
$newSerial = 45207

for ($row = 2; $row -le 101; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
